$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the character constellation value (column D, row 2) from 0 to 2
$ws.Range("D2").Value = 2

# Update the active selection to D2 (matches the diff's sheetView selection)
$ws.Range("D2").Select()
